$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.842.23'
$ws.Range("E2").Value = '  +10.57%  '

$ws.Range("D3").Value = '3.257.81'
$ws.Range("E3").Value = '  +5.97%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.556'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0962'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +13.19%  '

$ws.Range("E12").Value = '  +2.24%  '

$ws.Range("D13").Value = '3.766.13'
$ws.Range("E13").Value = '  +5.84%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.01%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.77%  '

$ws.Range("D16").Value = '3.253.40'
$ws.Range("E16").Value = '  +5.92%  '

$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.31%  '

$ws.Range("D19").Value = '56.669.65'
$ws.Range("E19").Value = '  +10.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000106'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '304.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +15.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.00%  '

$ws.Range("E25").Value = '  -0.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("E28").Value = '  +4.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.38%  '

$ws.Range("E30").Value = '  +3.48%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  +4.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0483'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.82%  '

$ws.Range("E38").Value = '  +24.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '134.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.77%  '

$ws.Range("E42").Value = '  +4.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.120'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.280'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.45%  '

$ws.Range("D48").Value = '2.151.61'
$ws.Range("E48").Value = '  +3.84%  '

$ws.Range("E49").Value = '  +2.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +43.83%  '

$ws.Range("E51").Value = '  -4.21%  '
